$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.192.37"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "1.857.58"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").Value = "238.40"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "0.6910"
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("D8").Value = "0.07741"
$ws.Range("E8").Value = "  +4.99%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("D11").Value = "0.08071"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "1.867.15"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "0.7235"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").Value = "5.205"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "89.56"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "29.198.85"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").Value = "5.751"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("D18").Value = "0.000007811"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").Value = "13.29"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").Value = "235.58"
$ws.Range("E20").Value = "  -2.93%  "
$ws.Range("D21").Value = "0.9991"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").Value = "2.105.02"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("D24").Value = "7.460"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("D25").Value = "161.91"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "8.988"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "0.1442"
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("D28").Value = "18.09"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "1.963"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("D30").Value = "1.401"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D31").Value = "4.531"
$ws.Range("E31").Value = "  +3.11%  "
$ws.Range("D32").Value = "1.488"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").Value = "0.05182"
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("D35").Value = "1.186"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").Value = "0.7056"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("D37").Value = "1.024"
$ws.Range("E37").Value = "  +3.36%  "
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "2.683"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").Value = "0.9245"
$ws.Range("E41").Value = "  +6.68%  "
$ws.Range("D42").Value = "1.097.49"
$ws.Range("E42").Value = "  +7.42%  "
$ws.Range("D43").Value = "5.986"
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D44").Value = "0.4292"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("D45").Value = "70.63"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("D46").Value = "0.9997"
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("D47").Value = "102.16"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "1.795"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("D49").Value = "2.001.94"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "9.182"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").Value = "7.010"
$ws.Range("E51").Value = "  -3.49%  "
